# Update the confidentiality footer date and the model's weight /
# percent-change figures on Sheet1. The sheet is protected, so we must
# unprotect, make the edits, then restore protection with the same
# password.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect("D382")

# Footer text: 2021-04-23 -> 2021-04-26
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-26 for illustrative purposes only and are subject to change."

# Row 2 (XLB / Materials Select Sector SPDR / Basic Materials Funds)
$ws.Range("D2").Value = 0.2521009190876301
$ws.Range("E2").Value = 0.005910022916415514

# Row 3 (XLF / Financial Select Sector SPDR Fund / Financial Services Funds)
$ws.Range("D3").Value = 0.2482789656632612
$ws.Range("E3").Value = 0.003957037874505298

# Row 4 (XLI / Industrial Select Sector SPDR Fund / Industrials Funds)
$ws.Range("D4").Value = 0.247575559570084
$ws.Range("E4").Value = -0.002953918865695071

# Row 5 (XLC / Communication Services Select Sector SPDR Fund / Telecommunication Funds)
$ws.Range("D5").Value = 0.2520445556790247
$ws.Range("E5").Value = 0.002092871157619358

# Row 6 (Total)
$ws.Range("E6").Value = 0.0022685501445463

$ws.Protect("D382")
